$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Tipo de documento" column (E) used the shared label "Cedula" for every
# record (E2:E4). The relabelling to "CC" touches every cell that held the
# old text so the underlying shared string is updated in place instead of
# leaving stray duplicates behind.
$ws.Range("E2").Value = "CC"
$ws.Range("E3").Value = "CC"
$ws.Range("E4").Value = "CC"

# Cursor/selection ends on D1 after the edit.
$ws.Range("D1").Select()
